# Link farms to pasture objects: add a new "Pasture" column (D) that
# classifies each farm as either "Natural Pasture" or "Sown Permanent Pasture".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("D1").Value = "Pasture"

# Fill in pasture classification for each farm row (2-43), alternating between
# "Natural Pasture" (even rows) and "Sown Permanent Pasture" (odd rows)
for ($r = 2; $r -le 43; $r++) {
    if ($r % 2 -eq 0) {
        $ws.Range("D$r").Value = "Natural Pasture"
    } else {
        $ws.Range("D$r").Value = "Sown Permanent Pasture"
    }
}

# Match the source workbook's column D width as closely as the engine allows
$ws.Columns.Item(4).ColumnWidth = 21.6

# Match the updated selection state (cell C2 selected)
$ws.Range("C2").Select()
